$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for the
#    2022-Q4 quarter just below the header row, pushing the
#    existing quarter rows down by one.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Re-apply the same formatting used by the other rows' first column
# (bold, bordered, centered) to the new A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 14
$summary.Range("D2").Value = 5.64

# ---------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" worksheet right after "总计" and
#    before the (currently second) "2022-Q3" sheet, by duplicating
#    the existing "2022-Q3" sheet (to inherit all formatting) and
#    then overwriting its contents with the 2022-Q4 fund data.
# ---------------------------------------------------------------
$existingQ3 = $wb.Worksheets.Item(2)
$existingQ3.Copy($existingQ3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The 2022-Q3 sheet had 16 data rows (rows 2-17); 2022-Q4 only has
# 14 data rows (rows 2-15), so remove the two extra rows.
$newSheet.Rows("16:17").Delete()

# Columns B, D, E, F and G hold text-formatted values (fund code,
# and numbers kept as formatted strings) - force them to Text so
# values such as "010094" or "1.60" are not coerced into numbers.
$newSheet.Range("B2:B15").NumberFormat = "@"
$newSheet.Range("D2:G15").NumberFormat = "@"

$newSheet.Cells.Item(2, 2).Value = "010094"
$newSheet.Cells.Item(2, 3).Value = "交银施罗德产业机遇混合"
$newSheet.Cells.Item(2, 4).Value = "16.65"
$newSheet.Cells.Item(2, 5).Value = "85.76"
$newSheet.Cells.Item(2, 6).Value = "7.46"
$newSheet.Cells.Item(2, 7).Value = "1.2421"
$newSheet.Cells.Item(2, 8).Value = 2
$newSheet.Cells.Item(3, 2).Value = "000021"
$newSheet.Cells.Item(3, 3).Value = "华夏优势增长混合"
$newSheet.Cells.Item(3, 4).Value = "53.36"
$newSheet.Cells.Item(3, 5).Value = "91.07"
$newSheet.Cells.Item(3, 6).Value = "1.98"
$newSheet.Cells.Item(3, 7).Value = "1.0565"
$newSheet.Cells.Item(3, 8).Value = 7
$newSheet.Cells.Item(4, 2).Value = "519773"
$newSheet.Cells.Item(4, 3).Value = "交银施罗德数据产业灵活配置混合A"
$newSheet.Cells.Item(4, 4).Value = "13.79"
$newSheet.Cells.Item(4, 5).Value = "86.66"
$newSheet.Cells.Item(4, 6).Value = "7.55"
$newSheet.Cells.Item(4, 7).Value = "1.0411"
$newSheet.Cells.Item(4, 8).Value = 2
$newSheet.Cells.Item(5, 2).Value = "519732"
$newSheet.Cells.Item(5, 3).Value = "交银定期支付双息平衡混合"
$newSheet.Cells.Item(5, 4).Value = "39.61"
$newSheet.Cells.Item(5, 5).Value = "69.68"
$newSheet.Cells.Item(5, 6).Value = "1.60"
$newSheet.Cells.Item(5, 7).Value = "0.6338"
$newSheet.Cells.Item(5, 8).Value = 10
$newSheet.Cells.Item(6, 2).Value = "010180"
$newSheet.Cells.Item(6, 3).Value = "华夏科技龙头两年定期开放混合"
$newSheet.Cells.Item(6, 4).Value = "14.70"
$newSheet.Cells.Item(6, 5).Value = "94.82"
$newSheet.Cells.Item(6, 6).Value = "3.29"
$newSheet.Cells.Item(6, 7).Value = "0.4836"
$newSheet.Cells.Item(6, 8).Value = 9
$newSheet.Cells.Item(7, 2).Value = "000061"
$newSheet.Cells.Item(7, 3).Value = "华夏盛世混合"
$newSheet.Cells.Item(7, 4).Value = "14.79"
$newSheet.Cells.Item(7, 5).Value = "80.71"
$newSheet.Cells.Item(7, 6).Value = "2.56"
$newSheet.Cells.Item(7, 7).Value = "0.3786"
$newSheet.Cells.Item(7, 8).Value = 3
$newSheet.Cells.Item(8, 2).Value = "012173"
$newSheet.Cells.Item(8, 3).Value = "国泰兴泽优选一年持有期混合A"
$newSheet.Cells.Item(8, 4).Value = "8.18"
$newSheet.Cells.Item(8, 5).Value = "92.65"
$newSheet.Cells.Item(8, 6).Value = "3.45"
$newSheet.Cells.Item(8, 7).Value = "0.2822"
$newSheet.Cells.Item(8, 8).Value = 8
$newSheet.Cells.Item(9, 2).Value = "012174"
$newSheet.Cells.Item(9, 3).Value = "国泰兴泽优选一年持有期混合C"
$newSheet.Cells.Item(9, 4).Value = "5.86"
$newSheet.Cells.Item(9, 5).Value = "92.65"
$newSheet.Cells.Item(9, 6).Value = "3.45"
$newSheet.Cells.Item(9, 7).Value = "0.2022"
$newSheet.Cells.Item(9, 8).Value = 8
$newSheet.Cells.Item(10, 2).Value = "014549"
$newSheet.Cells.Item(10, 3).Value = "交银施罗德数据产业灵活配置混合C"
$newSheet.Cells.Item(10, 4).Value = "2.43"
$newSheet.Cells.Item(10, 5).Value = "86.66"
$newSheet.Cells.Item(10, 6).Value = "7.55"
$newSheet.Cells.Item(10, 7).Value = "0.1835"
$newSheet.Cells.Item(10, 8).Value = 2
$newSheet.Cells.Item(11, 2).Value = "001924"
$newSheet.Cells.Item(11, 3).Value = "华夏国企改革灵活配置混合"
$newSheet.Cells.Item(11, 4).Value = "2.49"
$newSheet.Cells.Item(11, 5).Value = "85.98"
$newSheet.Cells.Item(11, 6).Value = "3.84"
$newSheet.Cells.Item(11, 7).Value = "0.0956"
$newSheet.Cells.Item(11, 8).Value = 3
$newSheet.Cells.Item(12, 2).Value = "002292"
$newSheet.Cells.Item(12, 3).Value = "诺安益鑫灵活配置混合A"
$newSheet.Cells.Item(12, 4).Value = "0.37"
$newSheet.Cells.Item(12, 5).Value = "69.58"
$newSheet.Cells.Item(12, 6).Value = "4.94"
$newSheet.Cells.Item(12, 7).Value = "0.0183"
$newSheet.Cells.Item(12, 8).Value = 3
$newSheet.Cells.Item(13, 2).Value = "008336"
$newSheet.Cells.Item(13, 3).Value = "宝盈祥裕增强回报混合A"
$newSheet.Cells.Item(13, 4).Value = "0.69"
$newSheet.Cells.Item(13, 5).Value = "36.55"
$newSheet.Cells.Item(13, 6).Value = "2.51"
$newSheet.Cells.Item(13, 7).Value = "0.0173"
$newSheet.Cells.Item(13, 8).Value = 8
$newSheet.Cells.Item(14, 2).Value = "008337"
$newSheet.Cells.Item(14, 3).Value = "宝盈祥裕增强回报混合C"
$newSheet.Cells.Item(14, 4).Value = "0.08"
$newSheet.Cells.Item(14, 5).Value = "36.55"
$newSheet.Cells.Item(14, 6).Value = "2.51"
$newSheet.Cells.Item(14, 7).Value = "0.0020"
$newSheet.Cells.Item(14, 8).Value = 8
$newSheet.Cells.Item(15, 2).Value = "014550"
$newSheet.Cells.Item(15, 3).Value = "诺安益鑫灵活配置混合C"
$newSheet.Cells.Item(15, 4).Value = "0.02"
$newSheet.Cells.Item(15, 5).Value = "69.58"
$newSheet.Cells.Item(15, 6).Value = "4.94"
$newSheet.Cells.Item(15, 7).Value = "0.0010"
$newSheet.Cells.Item(15, 8).Value = 3

# Restore the default (unstyled) cell style now that the text
# values are safely stored as strings, matching the look of the
# other quarter sheets.
$newSheet.Range("B2:B15").Style = "Normal"
$newSheet.Range("D2:G15").Style = "Normal"
